# Updated symbol list on Wed Dec 14 06:43:56 UTC 2022 with GitHub Actions
#
# The "Price" column (D) on the active sheet holds scraped coin prices as
# literal text (e.g. "274.58"), not numbers - trailing zeros and exact
# digit counts are meaningful. We refresh each changed quote in place,
# forcing the destination cell to Text format first so Excel/COM does not
# silently re-parse the digit string into a floating point number (which
# would drop significant trailing zeros such as 0.03080 -> 0.0308).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = @(
    @("D2",  "275.40"),
    @("D3",  "23.03"),
    @("D4",  "6.393"),
    @("D5",  "0.06267"),
    @("D6",  "3.658"),
    @("D7",  "6.666"),
    @("D8",  "1.362"),
    @("D9",  "0.8312"),
    @("D10", "0.01376"),
    @("D12", "0.08319"),
    @("D13", "0.03426"),
    @("D14", "0.03080"),
    @("D15", "0.09307"),
    @("D16", "3.858"),
    @("D17", "0.001642"),
    @("D18", "0.04767"),
    @("D19", "0.006316"),
    @("D21", "0.001093"),
    @("D23", "3.715"),
    @("D24", "2.353"),
    @("D40", "0.04696"),
    @("D41", "0.007077"),
    @("D43", "0.003699"),
    @("D44", "0.01211"),
    @("D45", "0.00006250"),
    @("D48", "0.7699"),
    @("D49", "0.04434")
)

foreach ($update in $priceUpdates) {
    $cellRef = $update[0]
    $newPrice = $update[1]

    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $newPrice
}
